# Apply the "configuration table function is completed" edit.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IteratValueInfo")

# Drop the stale row-outline grouping metadata left over on the sheet
# (no row actually carries a non-default outline level any more); this
# leaves the (still valid) column outline level untouched.
$ws.Outline.ShowLevels(0)

# The shared string previously used by row 8 ("BeInjuredIntervalTime")
# is renamed to "IntervalTime".
$ws.Range("A8").Value = "IntervalTime"

# Append the new config entry row.
$ws.Range("A9").Value = "LevelTime"
$ws.Range("B9").Value = 10
$ws.Range("C9").Value = 0

# Update selection to reflect the next empty row, as in the authored file.
$ws.Range("A10").Select()

# Match the saved window size recorded in the workbook view.
$win = $excel.ActiveWindow
$win.Width = 13070
$win.Height = 4410

